$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 163.33333
$ws.Range("I2").Value = 163.33333
$ws.Range("K2").Value = 163.33333
$ws.Range("M2").Value = -50.33332999999999

$ws.Range("H46").Value = 10200
$ws.Range("J46").Value = 8400
$ws.Range("L46").Value = 25200
$ws.Range("N46").Value = -25438

$ws.Range("H53").Value = 871.61536
$ws.Range("I53").Value = 893.7273
$ws.Range("K53").Value = 893.7273
$ws.Range("M53").Value = -256.7273

$ws.Range("H60").Value = 10200
$ws.Range("J60").Value = 8400
$ws.Range("L60").Value = 25200
$ws.Range("N60").Value = -26168

$ws.Range("H96").Value = 542
$ws.Range("I96").Value = 610.4
$ws.Range("J96").Value = 200
$ws.Range("K96").Value = 1831.2
$ws.Range("L96").Value = 600
$ws.Range("M96").Value = -458.1999999999998
$ws.Range("N96").Value = -3346

$ws.Range("H98").Value = 2187.25
$ws.Range("I98").Value = 2583
$ws.Range("K98").Value = 2583
$ws.Range("M98").Value = -1085

$ws.Range("H122").Value = 2187.25
$ws.Range("I122").Value = 2583
$ws.Range("K122").Value = 7749
$ws.Range("M122").Value = -5299

$ws.Range("H131").Value = 967.5
$ws.Range("I131").Value = 967.5
$ws.Range("K131").Value = 2902.5
$ws.Range("M131").Value = 2137.5

$ws.Range("H137").Value = 1724.75
$ws.Range("I137").Value = 1724.75
$ws.Range("K137").Value = 5174.25
$ws.Range("M137").Value = -2624.25

$ws.Range("H141").Value = 550
$ws.Range("I141").Value = 550
$ws.Range("K141").Value = 1650
$ws.Range("M141").Value = 3530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N61").ClearContents()
$ws.Range("H61").Value = 2999
$ws.Range("I61").Value = 2999
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2787

$ws.Range("H122").Value = 2870.6
$ws.Range("I122").Value = 2613.6667
$ws.Range("K122").Value = 7841.000100000001
$ws.Range("M122").Value = -5391.000100000001

$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 2999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8997
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13618.077
$ws.Range("I105").Value = 14436.25
$ws.Range("K105").Value = 14436.25
$ws.Range("M105").Value = -12689.25

$ws.Range("H134").Value = 4012
$ws.Range("I134").Value = 4012
$ws.Range("K134").Value = 12036
$ws.Range("M134").Value = -9501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 39500
$ws.Range("J74").Value = 39500
$ws.Range("L74").Value = 39500
$ws.Range("N74").Value = -41248

$ws.Range("H77").Value = 39500
$ws.Range("J77").Value = 39500
$ws.Range("L77").Value = 118500
$ws.Range("N77").Value = -127236

$ws.Range("H105").Value = 825
$ws.Range("I105").Value = 750
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 750
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = 997
$ws.Range("N105").Value = -4394

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M68").ClearContents()
$ws.Range("H68").Value = 1001.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1001.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3004.5
$ws.Range("N68").Value = -4626.5

$ws.Range("M71").ClearContents()
$ws.Range("H71").Value = 1001.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1001.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 9013.5
$ws.Range("N71").Value = -17125.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M18").ClearContents()
$ws.Range("H18").Value = 4503
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 4503
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 4503
$ws.Range("N18").Value = -5089

$ws.Range("H22").Value = 17669
$ws.Range("J22").Value = 14999
$ws.Range("L22").Value = 14999
$ws.Range("N22").Value = -16057

$ws.Range("H97").Value = 723.3333
$ws.Range("I97").Value = 723.3333
$ws.Range("K97").Value = 723.3333
$ws.Range("M97").Value = -227.3333

$ws.Range("H122").Value = 2503.5
$ws.Range("I122").Value = 2007
$ws.Range("K122").Value = 6021
$ws.Range("M122").Value = -3571

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 22800
$ws.Range("J4").Value = 22800
$ws.Range("L4").Value = 22800
$ws.Range("N4").Value = -23026

$ws.Range("N7").ClearContents()
$ws.Range("H7").Value = 2766.6667
$ws.Range("I7").Value = 2766.6667
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2766.6667
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2654.6667

$ws.Range("N16").ClearContents()
$ws.Range("H16").Value = 1649.8
$ws.Range("I16").Value = 1649.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1649.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1479.8

$ws.Range("H28").Value = 22800
$ws.Range("J28").Value = 22800
$ws.Range("L28").Value = 22800
$ws.Range("N28").Value = -23264

$ws.Range("H31").Value = 33999.5
$ws.Range("J31").Value = 33999.5
$ws.Range("L31").Value = 33999.5
$ws.Range("N31").Value = -34495.5

$ws.Range("H37").Value = 22800
$ws.Range("J37").Value = 22800
$ws.Range("L37").Value = 22800
$ws.Range("N37").Value = -23014

$ws.Range("N40").ClearContents()
$ws.Range("H40").Value = 7500
$ws.Range("I40").Value = 7500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7364

$ws.Range("M61").ClearContents()
$ws.Range("H61").Value = 5833.6665
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 5833.6665
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 5833.6665
$ws.Range("N61").Value = -6237.6665

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

$ws.Range("H64").Value = 63126
$ws.Range("J64").Value = 63126
$ws.Range("L64").Value = 63126
$ws.Range("N64").Value = -63576

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240

$ws.Range("H67").Value = 63126
$ws.Range("J67").Value = 63126
$ws.Range("L67").Value = 63126
$ws.Range("N67").Value = -64686

$ws.Range("H68").Value = 3500
$ws.Range("I68").Value = 2750
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2750
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -2001
$ws.Range("N68").Value = -5498

$ws.Range("H71").Value = 3500
$ws.Range("I71").Value = 2750
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 13750
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -10006
$ws.Range("N71").Value = -27488

$ws.Range("M113").ClearContents()
$ws.Range("H113").Value = 5833.6665
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5833.6665
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5833.6665
$ws.Range("N113").Value = -10173.6665

$ws.Range("N126").ClearContents()
$ws.Range("H126").Value = 2766.6667
$ws.Range("I126").Value = 2766.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8300.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5830.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 166672.67
$ws.Range("I28").Value = 149999
$ws.Range("K28").Value = 149999
$ws.Range("M28").Value = -149651

$ws.Range("H63").Value = 25000
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 25000
$ws.Range("N63").Value = -26248

$ws.Range("H66").Value = 25000
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 75000
$ws.Range("N66").Value = -81240

$ws.Range("H68").Value = 29635.5
$ws.Range("J68").Value = 29635.5
$ws.Range("L68").Value = 29635.5
$ws.Range("N68").Value = -31257.5

$ws.Range("H71").Value = 29635.5
$ws.Range("J71").Value = 29635.5
$ws.Range("L71").Value = 88906.5
$ws.Range("N71").Value = -97018.5

$ws.Range("H106").Value = 80000
$ws.Range("J106").Value = 80000
$ws.Range("L106").Value = 80000
$ws.Range("N106").Value = -82524

$ws.Range("H109").Value = 187777
$ws.Range("J109").Value = 187777
$ws.Range("L109").Value = 187777
$ws.Range("N109").Value = -190551

$ws.Range("H122").Value = 1635
$ws.Range("I122").Value = 452.5
$ws.Range("K122").Value = 1357.5
$ws.Range("M122").Value = 1092.5

$ws.Range("H126").Value = 5400
$ws.Range("I126").Value = 800
$ws.Range("K126").Value = 2400
$ws.Range("M126").Value = 70

$ws.Range("H136").Value = 2870.4285
$ws.Range("I136").Value = 1818.8
$ws.Range("J136").Value = 5499.5
$ws.Range("K136").Value = 5456.4
$ws.Range("L136").Value = 16498.5
$ws.Range("M136").Value = -2906.4
$ws.Range("N136").Value = -21598.5
